$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-decimal-looking price cells to stay text (matching the
# source feed's formatting) instead of Excel auto-converting them to
# numbers, which would also introduce float round-trip artifacts.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.586.04'
$ws.Range("E2").Value = '  -1.39%  '
$ws.Range("D3").Value = '1.666.49'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '215.29'
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").Value = '0.515'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '23.60'
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("D11").Value = '0.0883'
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("D12").Value = '1.901.84'
$ws.Range("E12").Value = '  -3.49%  '
$ws.Range("D13").Value = '1.692.19'
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("E14").Value = '  -2.59%  '
$ws.Range("D15").Value = '0.557'
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").Value = '66.26'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '250.06'
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D18").Value = '27.614.60'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("E19").Value = '  -3.52%  '
$ws.Range("D20").Value = '7.55'
$ws.Range("E20").Value = '  -4.30%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  -3.40%  '
$ws.Range("E23").Value = '  -4.76%  '
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -5.25%  '
$ws.Range("D25").Value = '146.52'
$ws.Range("E25").Value = '  -2.07%  '
$ws.Range("D26").Value = '16.49'
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -2.51%  '
$ws.Range("E30").Value = '  +3.84%  '
$ws.Range("D31").Value = '0.0508'
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("D33").Value = '1.472.28'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("E34").Value = '  -5.43%  '
$ws.Range("E35").Value = '  -5.27%  '
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("E38").Value = '  -6.12%  '
$ws.Range("D39").Value = '0.0172'
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("D40").Value = '69.63'
$ws.Range("E40").Value = '  -2.67%  '
$ws.Range("E41").Value = '  -5.62%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -7.29%  '
$ws.Range("D44").Value = '1.809.88'
$ws.Range("E44").Value = '  -3.45%  '
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("D46").Value = '0.789'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("E47").Value = '  -3.54%  '
$ws.Range("D48").Value = '89.42'
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("D49").Value = '0.0₆0108'
$ws.Range("E49").Value = '  -3.40%  '

# Rows 50 and 51 swap contents: Algorand moves up to rank 48 (row 50)
# and BitcoinSV drops to rank 49 (row 51), each with refreshed figures.
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -3.12%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = '41.55'
$ws.Range("E51").Value = '  +14.47%  '
